$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")

# Final target rows 2-52: Titre, Jours en Hausse, Jours en Baisse, Variation Totale (%), Derniere Variation (%), Recommandation, Strategie
$data = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3365.25, 108, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2770, 690, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2690, 685, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2581.24, 647.52, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI", 0, 2, 1880, 940, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 2, 1860, 935, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1646.67, 425.59, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1463.74, 369.61, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1362.53, 340.15, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 570.69, 142.38, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 547.54, 136.63, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 535.96, 133.99, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 529.71, 139.02, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 526.72, 131.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 430.31, 108.15, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 381.82, 95.06, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE            (**)", 0, 2, 371.92, 185.85, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                    (**)", 0, 1, 218.25, 218.25, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                 (**)", 0, 1, 216.31, 216.31, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                       (**)", 0, 1, 214.57, 214.57, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                         (**)", 0, 1, 214.08, 214.08, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                           (**)", 0, 1, 192.4, 192.4, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                          (**)", 0, 1, 192.26, 192.26, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 1, 192.05, 192.05, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                       (**)", 0, 1, 191.44, 191.44, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE             (**)", 0, 1, 188.06, 188.06, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE        (**)", 0, 1, 184.97, 184.97, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 4, 0, 28.2, 6.56, "🟢 Achat", "✅ Renforcer"),
    @("ORAGROUP TOGO (ORGT)", 2, 0, 14.95, 7.47, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 2, 0, 14.31, 7.07, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 2, 0, 13.19, 7.43, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 2, 0, 11.71, 6.49, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 2, 1, 10.11, 7.5, "🟡 Observer", "👀 À surveiller"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 1, 0, 5.88, 5.88, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 1, 1, 4.5, -2.96, "🟡 Observer", "👀 À surveiller"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 1, 0, 3.9, 3.9, "🟡 Observer", "➖ Neutre"),
    @("SAPH CI (SPHC)", 1, 1, 1.47, -1.72, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("PALM CI (PALC)", 0, 1, -1.1, -1.1, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 0, 1, -1.15, -1.15, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA ML (BOAM)", 1, 1, -1.39, 3.49, "🟡 Observer", "👀 À surveiller"),
    @("AIR LIQUIDE CI (SIVC)", 0, 1, -1.44, -1.44, "🟡 Observer", "➖ Neutre"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -2.03, -2.03, "🟡 Observer", "➖ Neutre"),
    @("CORIS BANK INTERNATIONAL (CBIBF)", 1, 1, -2.94, 4.52, "🟡 Observer", "👀 À surveiller"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 0, 1, -3.42, -3.42, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -3.7, -3.7, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 0, 2, -4.42, -2.6, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 0, 1, -4.78, -4.78, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 0, 1, -5.84, -5.84, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA BF (BOABF)", 0, 1, -7.42, -7.42, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 0, 3, -22.34, -7.49, "🔴 Vente", "⚠️ Risque de décrochage")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

$ws2 = $wb.Worksheets.Item("Top_YTD")

# Final target rows 2-11: Titre, Progression YTD (%)
$data2 = @(
    @("BRVM - SERVICES PUBLICS", 9125249.29),
    @("AIR LIQUIDE CI", 394346.8),
    @("NEI-CEDA CI", 355742.07),
    @("BRVM - AUTRES SECTEURS", 308428.46),
    @("BRVM - DISTRIBUTION", 68398.86),
    @("BRVM - TRANSPORT", 47027.17),
    @("BRVM - AGRICULTURE", 37594.43),
    @("CFAO MOTORS CI", 10716),
    @("SETAO CI", 10508.75),
    @("BRVM - INDUSTRIELS", 3367.99)
)

$r2 = 2
foreach ($row in $data2) {
    $ws2.Cells.Item($r2, 1).Value = $row[0]
    $ws2.Cells.Item($r2, 2).Value = $row[1]
    $r2 = $r2 + 1
}
